# Update cryptocurrency price/volume data to reflect the latest scrape (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.380.79"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "'1.574.58"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'290.50"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "'0.3761"
$ws.Range("E7").Value = "  +2.70%  "
$ws.Range("D8").Value = "'50.09"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "'0.3416"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'1.169"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").Value = "'0.07681"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "'21.37"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "'5.996"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'6.943"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "'1.572.86"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'90.38"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'0.06728"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "'16.76"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").Value = "'6.242"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'0.5281"
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("D24").Value = "'12.03"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "'22.381.53"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "'2.784"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "'20.26"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "'144.66"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'5.059"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").Value = "'126.50"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'1.748.50"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'1.036"
$ws.Range("E33").Value = "  +12.66%  "
$ws.Range("D34").Value = "'6.279"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "'2.024"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").Value = "'0.08542"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'0.02552"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "'0.2331"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06552"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "'5.522"
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("D42").Value = "'1.300"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").Value = "'11.69"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "'0.6453"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "'14.15"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'0.6041"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").Value = "'3.779"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'1.307"
$ws.Range("E49").Value = "  +11.24%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'125.52"
$ws.Range("E51").Value = "  +1.91%  "
